# Update the marksheet's correct/total marks figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# "Marking" row - correct answers count: 3 -> 5
$ws.Range("B11").Value = 5

# "Total" row - total score: 63 -> 105
$ws.Range("B12").Value = 105

# "Total" row - correct/total marks display: "61/84" -> "105/140"
$ws.Range("E12").Value = "105/140"
